$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "[-, 'MEC-3B-M.S.R. ar Cond.', 'MEC-3B-Mec. Manut.Equip. ind.', -]"

$ws.Range("E6").Value = "-"

$ws.Range("C7").Value = "['MEC-3B-Mec. Manut.Equip. ind.', 'MEC-3B-M.S.R. ar Cond.', -, -]"
$ws.Range("E7").Value = "-"

$ws.Range("B8").Value = "[-, -, 'MEC-3B-M.S.R. ar Cond.', 'MEC-3B-Mec. Manut.Equip. ind.']"
$ws.Range("C8").Value = "['MEC-3B-Mec. Manut.Equip. ind.', 'MEC-3B-M.S.R. ar Cond.', -, -]"

$ws.Range("B18").Value = "['MEC-2NB-M. Maq. E. I.', 'MEC-2NB-M.S.R.A.C.', 'ELM-2NA-Sistemas de Refrigeração', -]"
$ws.Range("C18").Value = "['MEC-2NA-M.S.R.A.C.', -, 'MEC-2NB-M.S.R.A.C.', 'ELM-2NA-Sistemas de Refrigeração']"
$ws.Range("D18").Value = "[-, -, 'MEC-2NA-M.S.R.A.C.', -]"

$ws.Range("B19").Value = "['MEC-2NB-M. Maq. E. I.', 'MEC-2NB-M.S.R.A.C.', -, -]"
$ws.Range("C19").Value = "[-, -, 'MEC-2NB-M. Maq. E. I.', -]"
$ws.Range("D19").Value = "[-, -, 'MEC-2NA-M.S.R.A.C.', -]"

$ws.Range("B20").Value = "['MEC-2NB-M. Maq. E. I.', 'MEC-2NB-M.S.R.A.C.', 'ELM-2NA-Sistemas de Refrigeração', -]"
$ws.Range("C20").Value = "[-, -, -, 'ELM-2NA-Sistemas de Refrigeração']"

$ws.Range("B21").Value = "[-, -, 'MEC-2NA-M.S.R.A.C.', -]"
